$d = $word.ActiveDocument

# 1 & 8: Title text appears twice (H1 heading and bold paragraph near the end) - identical change both places
$d.Content.Find.Execute("Play Don Corlimone for Free - Unique Fruit and Mafia Themed Slot", $true, $false, $false, $false, $false, $true, 1, $false, "Play Don Corlimone Free: Exciting Fruit-Mafia Slot Game", 2)

# 2: "What we like" bullet 1
$d.Content.Find.Execute("High-quality graphics and immersive sound effects", $true, $false, $false, $false, $false, $true, 1, $false, "High-quality graphics and attention to detail in symbols", 2)

# 3: "What we like" bullet 2
$d.Content.Find.Execute("Chance to trigger free spins and bonus rounds", $true, $false, $false, $false, $false, $true, 1, $false, "Catchy background music that creates atmosphere", 2)

# 4: "What we like" bullet 3
$d.Content.Find.Execute("Moderately high winning potential with 25 pay lines", $true, $false, $false, $false, $false, $true, 1, $false, "Moderately high winning potential with free spins and bonus rounds", 2)

# 5: "What we like" bullet 4
$d.Content.Find.Execute("Unique combination of fruit and mafia theme", $true, $false, $false, $false, $false, $true, 1, $false, "Unique combination of fruit and mafia themes", 2)

# 6: "What we don't like" bullet 1
$d.Content.Find.Execute("Custom symbols may initially cause confusion", $true, $false, $false, $false, $false, $true, 1, $false, "Initial confusion due to symbol customization", 2)

# 7: "What we don't like" bullet 2
$d.Content.Find.Execute("Background music can become overwhelming for some players", $true, $false, $false, $false, $false, $true, 1, $false, "Music can become overwhelming and needs to be manually turned off", 2)

# 9: meta description (italic paragraph)
$d.Content.Find.Execute("Read our review of Don Corlimone, the fruit and mafia themed slot game. Play for free and win big with its high-quality graphics and bonus features.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Don Corlimone, a unique fruit-mafia themed slot game. Play for free and enjoy exciting gameplay features.", 2)
